$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update company name (RAZON SOCIAL) text
$ws.Range("E7").Value = "DC DORAL CARTAGENA S.A.S."

# 2) Reverse the order of the "Periodo Mora" labels in E16:E35
#    Before: 2304, 2303, ... 2109 (descending)
#    After:  2109, 2110, ... 2304 (ascending)
$periodos = @(2109, 2110, 2111, 2112, 2201, 2202, 2203, 2204, 2205, 2206, 2207, 2208, 2209, 2210, 2211, 2212, 2301, 2302, 2303, 2304)
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = [string]$periodos[$i]
}

# 3) The special lower payment value (36667) that was tied to period 2304 (row 16)
#    now follows that period to its new row (35); row 16 reverts to the standard 100000
$ws.Range("F16").Value = 100000
$ws.Range("F35").Value = 36667

# 4) Update the mora value total on row 36 (second worker, period 2304)
$ws.Range("G36").Value = 1423500
